# Generate Report for Handback
# The "ea31f9a5-27b9-4331-94c3-1c02177394ae.md" file has now been handed
# back (it was previously only "Ready for handoff"), and the handback
# timestamps for both locales move forward to reflect the new handback.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet: both locale status columns for that file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# --- zh-cn sheet: status + handback datetime (column G = "Latest Handback DateTime") ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $handedBack
$zhcn.Range("G2").Value = "2016-03-02 15:39:47"
$zhcn.Range("G3").Value = "2016-03-02 15:39:47"

# --- de-de sheet: status + handback datetime (column G = "Latest Handback DateTime") ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $handedBack
$dede.Range("G2").Value = "2016-03-02 15:40:07"
$dede.Range("G3").Value = "2016-03-02 15:40:07"
